$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C6 (time) -> D6 formula (C6-B6) will recalc automatically
$ws.Range("C6").Value() = 0.60416666666666663

# Add new note text in E6, referencing (creating) a new shared string
$ws.Range("E6").Value() = "reacquainting myself with min organization, setting up any connect for max"

# Add a new row 7 with a time value in B7, matching the time format used in B6/C6
$ws.Range("B7").Value() = 0.70138888888888884
$ws.Range("B7").NumberFormat = $ws.Range("B6").NumberFormat

# Move the active selection to B8
$ws.Range("B8").Select()
